$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 15747
$ws1.Range("F9").Value = 15458
$ws1.Range("F11").Value = 9062
$ws1.Range("F18").Value = 203
$ws1.Range("F39").Value = 5570

# Sheet "全部类型"
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F5").Value = 15747
$ws2.Range("F9").Value = 15458
$ws2.Range("F11").Value = 9062
$ws2.Range("F18").Value = 203
$ws2.Range("F41").Value = 5570
